$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells already carrying each of the three named cell styles used on
# this sheet (plain / yellow fill / blue fill, all with the Meiryo font), so the
# new CJ cells can reuse the existing style entries instead of creating new ones.
$styleTemplates = @{ 1 = "A2"; 2 = "D2"; 3 = "N2" }

# New date column CJ for 2024/12/05.
$ws.Range("CJ1").ColumnWidth = 11.17
$ws.Range("CJ1").Value = "'2024/12/05"  # leading apostrophe forces text, like the other date headers
$ws.Range("CI1").Copy()
$ws.Range("CJ1").PasteSpecial(-4122)  # xlPasteFormats

# New daily values for 2024/12/05, one per row, with the style index that
# matches the original sheet (1 = plain, 2 = yellow fill, 3 = blue fill).
$rows = @(
    [PSCustomObject]@{ Row = 2; Value = 134.8; Style = 3 },
    [PSCustomObject]@{ Row = 3; Value = 217.6; Style = 1 },
    [PSCustomObject]@{ Row = 4; Value = 136.9; Style = 3 },
    [PSCustomObject]@{ Row = 5; Value = 143.3; Style = 1 },
    [PSCustomObject]@{ Row = 6; Value = 113.5; Style = 2 },
    [PSCustomObject]@{ Row = 7; Value = 183.1; Style = 1 },
    [PSCustomObject]@{ Row = 8; Value = 109.6; Style = 2 },
    [PSCustomObject]@{ Row = 9; Value = 212.5; Style = 1 },
    [PSCustomObject]@{ Row = 10; Value = 151; Style = 1 },
    [PSCustomObject]@{ Row = 11; Value = 123.3; Style = 2 },
    [PSCustomObject]@{ Row = 12; Value = 119; Style = 2 },
    [PSCustomObject]@{ Row = 13; Value = 127.3; Style = 3 },
    [PSCustomObject]@{ Row = 14; Value = 124.9; Style = 2 },
    [PSCustomObject]@{ Row = 15; Value = 161.9; Style = 1 },
    [PSCustomObject]@{ Row = 16; Value = 173; Style = 1 },
    [PSCustomObject]@{ Row = 17; Value = 141; Style = 1 },
    [PSCustomObject]@{ Row = 18; Value = 154.9; Style = 1 },
    [PSCustomObject]@{ Row = 19; Value = 161.5; Style = 1 },
    [PSCustomObject]@{ Row = 20; Value = 148.6; Style = 1 },
    [PSCustomObject]@{ Row = 21; Value = 141.8; Style = 1 },
    [PSCustomObject]@{ Row = 22; Value = 154.7; Style = 1 },
    [PSCustomObject]@{ Row = 23; Value = 188.7; Style = 1 },
    [PSCustomObject]@{ Row = 24; Value = 114.8; Style = 2 },
    [PSCustomObject]@{ Row = 25; Value = 163.4; Style = 1 },
    [PSCustomObject]@{ Row = 26; Value = 146.3; Style = 1 },
    [PSCustomObject]@{ Row = 27; Value = 138.8; Style = 3 },
    [PSCustomObject]@{ Row = 28; Value = 151.7; Style = 1 },
    [PSCustomObject]@{ Row = 29; Value = 150.1; Style = 1 },
    [PSCustomObject]@{ Row = 30; Value = 151.3; Style = 1 },
    [PSCustomObject]@{ Row = 31; Value = 140.3; Style = 1 },
    [PSCustomObject]@{ Row = 32; Value = 122.3; Style = 2 },
    [PSCustomObject]@{ Row = 33; Value = 143.4; Style = 1 },
    [PSCustomObject]@{ Row = 34; Value = 155.9; Style = 1 },
    [PSCustomObject]@{ Row = 35; Value = 160.4; Style = 1 },
    [PSCustomObject]@{ Row = 36; Value = 118.8; Style = 2 },
    [PSCustomObject]@{ Row = 37; Value = 178.7; Style = 1 },
    [PSCustomObject]@{ Row = 38; Value = 130.5; Style = 3 },
    [PSCustomObject]@{ Row = 39; Value = 162.7; Style = 1 },
    [PSCustomObject]@{ Row = 40; Value = 150.9; Style = 1 },
    [PSCustomObject]@{ Row = 41; Value = 134.4; Style = 3 },
    [PSCustomObject]@{ Row = 42; Value = 171.1; Style = 1 },
    [PSCustomObject]@{ Row = 43; Value = 156; Style = 1 },
    [PSCustomObject]@{ Row = 44; Value = 128.1; Style = 3 },
    [PSCustomObject]@{ Row = 45; Value = 170.6; Style = 1 },
    [PSCustomObject]@{ Row = 46; Value = 207.1; Style = 1 },
    [PSCustomObject]@{ Row = 47; Value = 158.9; Style = 1 },
    [PSCustomObject]@{ Row = 48; Value = 148; Style = 1 },
    [PSCustomObject]@{ Row = 49; Value = 152.2; Style = 1 },
    [PSCustomObject]@{ Row = 50; Value = 169.7; Style = 1 },
    [PSCustomObject]@{ Row = 51; Value = 120.3; Style = 2 },
    [PSCustomObject]@{ Row = 52; Value = 144; Style = 1 },
    [PSCustomObject]@{ Row = 53; Value = 112.4; Style = 2 }
)

foreach ($entry in $rows) {
    $dest = $ws.Range("CJ" + $entry.Row)
    $ws.Range($styleTemplates[$entry.Style]).Copy()
    $dest.PasteSpecial(-4122)  # xlPasteFormats
    $dest.Value = $entry.Value
}
